# Fix TTD buffer bug in ivtnmr.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the TTD label in B6 and its formula to subtract ProteinVolume
$ws.Range("B6").Value = "20x TTD pH 7.7 (minus protein volume)"
$ws.Range("G6").Formula = "=(TotalVolume-ProteinVolume)/F6"

# Update the "final" label in B5 and its row height
$ws.Range("B5").Value = "Final: Plasmd~30nM or Oligo~1uM."
$ws.Rows.Item(5).RowHeight = 40

# Move the active selection to D24 to match the saved view state
$ws.Range("D24").Select()
